$d = $word.ActiveDocument

function Replace-Exact($oldText, $newText) {
    $r = $d.Content
    $found = $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor text not found: $oldText"
    }
    $r.Text = $newText
}

function Replace-Span($firstText, $lastText, $newText) {
    $r1 = $d.Content
    $found1 = $r1.Find.Execute($firstText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found1) {
        throw "Span-start anchor not found: $firstText"
    }
    $start = $r1.Start

    $r2 = $d.Content
    $found2 = $r2.Find.Execute($lastText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found2) {
        throw "Span-end anchor not found: $lastText"
    }
    $end = $r2.End

    $merged = $d.Range($start, $end)
    $merged.Text = $newText
}

# ---- Title ----
Replace-Exact "Echoes of Eternity - A Celestial Symphony" "The Profound Insights of Mathematics and Its Widespread Impact"

# ---- Author name: "Jaime Torres" -> "Dr" + "." + " Albert Clayton" ----
$r = $d.Content
$r.Find.Execute("Jaime Torres", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Text = "Dr"
$pos = $r.End
$r2 = $d.Range($pos, $pos)
$r2.InsertAfter(".")
$pos2 = $r2.End
$r3 = $d.Range($pos2, $pos2)
$r3.InsertAfter(" Albert Clayton")

# ---- Email: "jm" / "torres@umich" / "edu" -> "ac" / "claytonphd@protonmail" / "com" ----
Replace-Exact "jm" "ac"
Replace-Exact "torres@umich" "claytonphd@protonmail"
Replace-Exact "edu" "com"

# ---- Body paragraph 1 ----
Replace-Exact "Drawn by a cosmic magnetism, we are entranced by the silent yet resounding melodies of the universe" "Mathematics, the language of the universe, unravels the intricate patterns and symmetries that permeate existence"

Replace-Exact " The very fabric of our universe vibrates with energy, orchestrating a Symphony of the Cosmos" " It is a boundless realm of exploration, where abstract concepts converge with tangible applications, illuminating the world around us"

Replace-Span " Throughout history, across cultures, humans have been bewitched by the celestial spectacle" " The history of space exploration can be seen as an effort to conduct this symphony, to understand and contribute to the grand cosmic opera, wherein celestial bodies play their individual musical parts" " From the cosmos' vast expanses to the intricate designs of nature, mathematics provides a framework for comprehending and harnessing the universe's underlying forces"

Replace-Exact "The celestial ambiance is woven with gravitational lullabies, where galaxies and stars dance in intricate waltzes" "Immersed in a world governed by numbers and equations, we discern the rhythmic harmonies of mathematical principles echoing throughout our lives"

Replace-Exact " The frequency of light from distant stars, the sonata of solar flares, the radioactive hum, and the melodic murmur of black holes - all these sounds, if transposed to human ears, would form a cacophony" " The Pythagorean theorem unveils the beauty of geometric relationships, guiding architects and engineers in constructing awe-inspiring structures"

Replace-Span " Yet, this vibrant chaos has structure, order, and rhythm, waiting to be analyzed and decoded" " As we probe deeper into the mysteries of the universe, we are becoming attuned to its musicality" " Calculus, a symphony of change, empowers scientists to model complex phenomena, opening doors to novel technological advancements"

Replace-Exact "On Earth, biology offers a terrestrial echo of the universe's Symphony" "Mathematics serves as a venerable instrument of discovery, propelling humanity's quest for knowledge"

Replace-Exact " Human cells exhibit rhythmic metabolic processes mimicking the pulsation of stars, and life's intricate system of interactions between organisms mirrors the celestial mechanics of planetary configurations" " It unveils the secrets of the cosmos, unraveling the mysteries of celestial bodies and guiding astronauts through the vast expanse of space"

Replace-Exact " By comprehending Earth's biological harmony, we deepen our understanding of the cosmic concord" " It unlocks the enigmas of subatomic particles, empowering physicists to explore the fundamental building blocks of matter"

Replace-Exact " Each species, like a unique instrument, adds its distinct timbre to the terrestrial symphony, contributing a unique rhythm to the symphony of life" " Mathematics reveals the intricate machinery of life, enabling biologists to decipher the genetic code and unravel the complexities of the human body"

# ---- Summary paragraph ----
Replace-Exact "In essence, the universe reverberates with an unspoken melody, an intrinsic cohesion between celestial bodies, energy fields, and life itself" "Mathematics, an intellectual odyssey, unveils the universe's mysteries, propels technological advancements, and underpins our understanding of the cosmos, nature, and life itself"

Replace-Span " Our goal is to capture this celestial symphony, to understand the rhythm of the universe, and to recognize our own place within its vast " " By listening to the murmurs of cosmos, we come closer to comprehending the fabric of the universe and our own role within it" " It remains an instrumental force in shaping our world, an enduring testament to the power of human intellect"

# ---- Add trailing empty paragraph ----
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endRange = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$endRange.InsertParagraphAfter()

Write-Output "edit complete"
